# Append a new row (row 42) to each of the 4 log sheets, replicating the
# previous last row (row 41) verbatim except for a new timestamp in column A.
# This mirrors a new sample having been appended to the logged data.

$wb = $excel.ActiveWorkbook

# --- ROW35-FE-LIFTER -------------------------------------------------
$ws = $wb.Worksheets.Item("ROW35-FE-LIFTER")
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A42").Value = 45745.32949260416
$ws.Range("B42").Value = "0x01,0x90"
$ws.Range("C42").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
$ws.Range("D42").Value = "0x01,0x7a"
$ws.Range("E42").Value = "0xd"
$ws.Range("F42").Value = 400
$ws.Range("G42").Value = 568631262647114.0 * 1000000000.0
$ws.Range("H42").Value = 378
$ws.Range("I42").Value = 13

# --- ROW35-MID-LIFTER --------------------------------------------------
$ws = $wb.Worksheets.Item("ROW35-MID-LIFTER")
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A42").Value = 45745.1770559838
$ws.Range("B42").Value = "0x01,0x90"
$ws.Range("C42").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws.Range("D42").Value = "0x01,0x7a"
$ws.Range("E42").Value = "0xe"
$ws.Range("F42").Value = 400
$ws.Range("G42").Value = 568631262647114.0 * 1000000000.0
$ws.Range("H42").Value = 378
$ws.Range("I42").Value = 14

# --- ROW02-FE-LIFTER -----------------------------------------------
$ws = $wb.Worksheets.Item("ROW02-FE-LIFTER")
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A42").Value = 45745.31992869213
$ws.Range("B42").Value = "0x01,0x90"
$ws.Range("C42").Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Range("D42").Value = "0x01,0x7a"
$ws.Range("E42").Value = "0x3"
$ws.Range("F42").Value = 400
$ws.Range("G42").Value = 568631262647114.0 * 1000000000.0
$ws.Range("H42").Value = 378
$ws.Range("I42").Value = 3

# --- ROW02-MID-LIFTER ----------------------------------------------
$ws = $wb.Worksheets.Item("ROW02-MID-LIFTER")
$ws.Range("A42").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A42").Value = 45745.3781578588
$ws.Range("B42").Value = "0x01,0x90"
$ws.Range("C42").Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Range("D42").Value = "0x01,0x7a"
$ws.Range("E42").Value = "0x3"
$ws.Range("F42").Value = 400
$ws.Range("G42").Value = 985046333984776.0 * 1000000000.0
$ws.Range("H42").Value = 378
$ws.Range("I42").Value = 3
